$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.492.03"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'1.836.80"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'261.96"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5385"
$ws.Range("E7").Value = "  +2.18%  "
$ws.Range("D8").Value = "'0.3002"
$ws.Range("E8").Value = "  -7.46%  "
$ws.Range("D9").Value = "'0.06929"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").Value = "'17.58"
$ws.Range("E10").Value = "  -7.30%  "
$ws.Range("D11").Value = "'0.7347"
$ws.Range("E11").Value = "  -6.28%  "
$ws.Range("D12").Value = "'1.839.63"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'0.07201"
$ws.Range("E13").Value = "  -7.21%  "
$ws.Range("D14").Value = "'89.20"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'13.83"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'0.000007902"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'26.521.87"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'2.074.87"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").Value = "'4.576"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "'5.989"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'9.208"
$ws.Range("D25").Value = "'142.89"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'2.174"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'1.715"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "'16.99"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'111.09"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'4.239"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'0.08856"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").Value = "'4.048"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'0.04842"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'2.931"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "'0.7289"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'3.093"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "'2.300"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "'0.01716"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("D40").Value = "'0.4716"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'0.9075"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'107.95"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").Value = "'5.908"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'7.422"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").Value = "'0.1252"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'9.026"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").Value = "'0.4081"
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").Value = "'34.82"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'0.8937"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'0.05767"
$ws.Range("E51").Value = "  -2.04%  "
